$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target range keeps its values as text (the source data stores
# prices/percentages as strings, e.g. "95.278.54" or "0.560"). Setting the
# NumberFormat to Text ("@") before assigning prevents Excel from silently
# re-interpreting these strings as numbers and losing formatting (trailing
# zeroes, multi-dot "thousand separators", etc). ClearFormats afterwards
# removes the temporary format so the cell style matches the original (no
# explicit style index), while the value stays text.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "95.278.54"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "3.579.91"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "236.80"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "651.05"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "0.400"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D10").Value = "1.01"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").Value = "3.579.67"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "42.40"
$ws.Range("E13").Value = "  -1.75%  "
$ws.Range("D14").Value = "6.44"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "4.244.13"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "94.976.46"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "3.567.53"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "7.93"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").Value = "12.53"
$ws.Range("E20").Value = "  -5.50%  "
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").Value = "3.47"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "508.55"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("D25").Value = "6.81"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("D27").Value = "95.33"
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("D28").Value = "12.58"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").Value = "3.771.94"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "3.02"
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("D32").Value = "11.46"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").Value = "31.84"
$ws.Range("E36").Value = "  +4.66%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "1.67"
$ws.Range("E37").Value = "  +11.97%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.560"
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("D39").Value = "8.55"
$ws.Range("E39").Value = "  +8.16%  "
$ws.Range("D40").Value = "582.43"
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").Value = "0.905"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").Value = "1.79"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "5.75"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").Value = "2.31"
$ws.Range("E46").Value = "  +5.71%  "
$ws.Range("D47").Value = "33.98"
$ws.Range("E47").Value = "  +30.49%  "
$ws.Range("D48").Value = "23.40"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("D49").Value = "0.0415"
$ws.Range("E49").Value = "  -4.24%  "
$ws.Range("D50").Value = "3.54"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "53.27"
$ws.Range("E51").Value = "  -0.93%  "

$dataRange.ClearFormats()
